$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.046916387279797
$ws.Range("D2").Value = 1.049188941271164
$ws.Range("E2").Value = 1.044417284361406
$ws.Range("I2").Value = 1.045312376455398
$ws.Range("J2").Value = 1.051968233168832
$ws.Range("K2").Value = 1.051946693467375
$ws.Range("L2").Value = 1.047188393011649
$ws.Range("N2").Value = 1.053462148722375

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.048096346228673
$ws.Range("D3").Value = 1.050117321147743
$ws.Range("E3").Value = 1.04543235793806
$ws.Range("I3").Value = 1.045706705812828
$ws.Range("J3").Value = 1.052795319690578
$ws.Range("K3").Value = 1.052686969253199
$ws.Range("L3").Value = 1.048014168458733
$ws.Range("N3").Value = 1.05429040980185

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.048859269755589
$ws.Range("D4").Value = 1.050717447142533
$ws.Range("E4").Value = 1.046088855880957
$ws.Range("I4").Value = 1.045960212026384
$ws.Range("J4").Value = 1.053329378647749
$ws.Range("K4").Value = 1.053164753955934
$ws.Range("L4").Value = 1.048547571809431
$ws.Range("N4").Value = 1.054825227183997

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.049179863834313
$ws.Range("D5").Value = 1.050969597854218
$ws.Range("E5").Value = 1.046364771524173
$ws.Range("I5").Value = 1.046066391334883
$ws.Range("J5").Value = 1.053553629730595
$ws.Range("K5").Value = 1.053365322730137
$ws.Range("L5").Value = 1.048771593370726
$ws.Range("N5").Value = 1.055049796729102

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.049233684899248
$ws.Range("D6").Value = 1.051011926751517
$ws.Range("E6").Value = 1.046411094545305
$ws.Range("I6").Value = 1.046084196180017
$ws.Range("J6").Value = 1.053591266832612
$ws.Range("K6").Value = 1.053398982018768
$ws.Range("L6").Value = 1.048809194636744
$ws.Range("N6").Value = 1.055087487280121

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.048863554092938
$ws.Range("D7").Value = 1.050720816950574
$ws.Range("E7").Value = 1.046092542976873
$ws.Range("I7").Value = 1.045961632349663
$ws.Range("J7").Value = 1.053332376150588
$ws.Range("K7").Value = 1.053167435111587
$ws.Range("L7").Value = 1.048550566064325
$ws.Range("N7").Value = 1.054828228943634

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.047315283049164
$ws.Range("D8").Value = 1.049502815916767
$ws.Range("E8").Value = 1.044760400690652
$ws.Range("I8").Value = 1.045445984177344
$ws.Range("J8").Value = 1.052247983800745
$ws.Range("K8").Value = 1.052197126974428
$ws.Range("L8").Value = 1.047467660620645
$ws.Range("N8").Value = 1.053742296632273

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.044582432421752
$ws.Range("D9").Value = 1.047351923730888
$ws.Range("E9").Value = 1.042410474658602
$ws.Range("I9").Value = 1.044524671996638
$ws.Range("J9").Value = 1.050328506074231
$ws.Range("K9").Value = 1.050477909390457
$ws.Range("L9").Value = 1.045552281478167
$ws.Range("N9").Value = 1.05182009302722

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.042757311877228
$ws.Range("D10").Value = 1.045914831617628
$ws.Range("E10").Value = 1.040842088437882
$ws.Range("I10").Value = 1.043901902190833
$ws.Range("J10").Value = 1.049042970396319
$ws.Range("K10").Value = 1.049325375710372
$ws.Range("L10").Value = 1.044270475540751
$ws.Range("N10").Value = 1.050532731741193

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.041966221505792
$ws.Range("D11").Value = 1.045291790456245
$ws.Range("E11").Value = 1.040162523076452
$ws.Range("I11").Value = 1.043630195791059
$ws.Range("J11").Value = 1.048484906242944
$ws.Range("K11").Value = 1.048824786171443
$ws.Range("L11").Value = 1.043714263749633
$ws.Range("N11").Value = 1.049973875072709

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.041672252331425
$ws.Range("D12").Value = 1.045060248113455
$ws.Range("E12").Value = 1.039910033981728
$ws.Range("I12").Value = 1.043528964109461
$ws.Range("J12").Value = 1.048277401617213
$ws.Range("K12").Value = 1.048638612989452
$ws.Range("L12").Value = 1.043507482908065
$ws.Range("N12").Value = 1.049766075766608

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.041735315324925
$ws.Range("D13").Value = 1.045109920052814
$ws.Range("E13").Value = 1.039964196864231
$ws.Range("I13").Value = 1.043550692593036
$ws.Range("J13").Value = 1.048321921807622
$ws.Range("K13").Value = 1.048678558286347
$ws.Range("L13").Value = 1.04355184622606
$ws.Range("N13").Value = 1.049810659180793

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.041941924446255
$ws.Range("D14").Value = 1.045272653490679
$ws.Range("E14").Value = 1.040141653650855
$ws.Range("I14").Value = 1.043621834233445
$ws.Range("J14").Value = 1.048467758236018
$ws.Range("K14").Value = 1.04880940177742
$ws.Range("L14").Value = 1.043697174841531
$ws.Range("N14").Value = 1.049956702713647

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.042069206769174
$ws.Range("D15").Value = 1.045372903374428
$ws.Range("E15").Value = 1.040250981540894
$ws.Range("I15").Value = 1.043665626129175
$ws.Range("J15").Value = 1.048557584368004
$ws.Range("K15").Value = 1.048889987996203
$ws.Range("L15").Value = 1.043786692830593
$ws.Range("N15").Value = 1.050046656409043

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.042809796836917
$ws.Range("D16").Value = 1.045956164456106
$ws.Range("E16").Value = 1.040887179467448
$ws.Range("I16").Value = 1.043919891330561
$ws.Range("J16").Value = 1.049079977248626
$ws.Range("K16").Value = 1.049358565737295
$ws.Range("L16").Value = 1.04430736446503
$ws.Range("N16").Value = 1.050569791147475

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.043274133082852
$ws.Range("D17").Value = 1.046321821436319
$ws.Range("E17").Value = 1.041286129686416
$ws.Range("I17").Value = 1.044078837583646
$ws.Range("J17").Value = 1.049407279538318
$ws.Range("K17").Value = 1.049652080304361
$ws.Range("L17").Value = 1.044633650661988
$ws.Range("N17").Value = 1.050897558243959

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.043544895440618
$ws.Range("D18").Value = 1.046535028686352
$ws.Range("E18").Value = 1.041518787890169
$ws.Range("I18").Value = 1.044171351188778
$ws.Range("J18").Value = 1.049598052650107
$ws.Range("K18").Value = 1.049823134386895
$ws.Range("L18").Value = 1.044823854000785
$ws.Range("N18").Value = 1.051088602275449

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.043637205405675
$ws.Range("D19").Value = 1.046607714288736
$ws.Range("E19").Value = 1.041598111128498
$ws.Range("I19").Value = 1.044202862540146
$ws.Range("J19").Value = 1.0496630781996
$ws.Range("K19").Value = 1.049881434345663
$ws.Range("L19").Value = 1.044888689153896
$ws.Range("N19").Value = 1.051153720168681

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.043224322195026
$ws.Range("D20").Value = 1.04628259759641
$ws.Range("E20").Value = 1.041243330533629
$ws.Range("I20").Value = 1.044061804544283
$ws.Range("J20").Value = 1.049372177260097
$ws.Range("K20").Value = 1.049620604287828
$ws.Range("L20").Value = 1.044598655032032
$ws.Range("N20").Value = 1.050862406116478

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.041881086625802
$ws.Range("D21").Value = 1.045224735792766
$ws.Range("E21").Value = 1.040089398947347
$ws.Range("I21").Value = 1.043600893298241
$ws.Range("J21").Value = 1.048424818983844
$ws.Range("K21").Value = 1.048770878043961
$ws.Range("L21").Value = 1.043654384140563
$ws.Range("N21").Value = 1.049913702482809

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.04103582823816
$ws.Range("D22").Value = 1.044558937845703
$ws.Range("E22").Value = 1.03936348062564
$ws.Range("I22").Value = 1.04330931866628
$ws.Range("J22").Value = 1.047827933950008
$ws.Range("K22").Value = 1.048235278833178
$ws.Range("L22").Value = 1.043059647053639
$ws.Range("N22").Value = 1.049315969803761

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.041483984153961
$ws.Range("D23").Value = 1.044911954665093
$ws.Range("E23").Value = 1.039748341706574
$ws.Range("I23").Value = 1.043464057055191
$ws.Range("J23").Value = 1.048144472386899
$ws.Range("K23").Value = 1.048519337822411
$ws.Range("L23").Value = 1.043375027197302
$ws.Range("N23").Value = 1.04963295776154

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.043246829829434
$ws.Range("D24").Value = 1.0463003213909
$ws.Range("E24").Value = 1.04126266976055
$ws.Range("I24").Value = 1.044069501650486
$ws.Range("J24").Value = 1.049388038891322
$ws.Range("K24").Value = 1.049634827402173
$ws.Range("L24").Value = 1.044614468402757
$ws.Range("N24").Value = 1.050878290273038

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.045289497983581
$ws.Range("D25").Value = 1.047908534263703
$ws.Range("E25").Value = 1.043018293206018
$ws.Range("I25").Value = 1.0447643593684
$ws.Range("J25").Value = 1.05082576862994
$ws.Range("K25").Value = 1.049652080304361
$ws.Range("L25").Value = 1.04604830891084
$ws.Range("N25").Value = 1.052318061752795
